$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 901.63635
$ws.Range("I53").Value = 935.3333
$ws.Range("K53").Value = 935.3333
$ws.Range("M53").Value = -298.3333
# Row 70
$ws.Range("H70").Value = 1554.1666
$ws.Range("J70").Value = 1500
$ws.Range("L70").Value = 4500
$ws.Range("N70").Value = -5040
# Row 73
$ws.Range("H73").Value = 1554.1666
$ws.Range("J73").Value = 1500
$ws.Range("L73").Value = 4500
$ws.Range("N73").Value = -6372
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 132
$ws.Range("H132").Value = 5172.75
$ws.Range("I132").Value = 5822.4287
$ws.Range("K132").Value = 17467.2861
$ws.Range("M132").Value = -14937.2861

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 109
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
# Row 110
$ws.Range("H110").Value = 511
$ws.Range("I110").Value = 511
$ws.Range("K110").Value = 511
$ws.Range("M110").Value = 1534

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 117
$ws.Range("H117").Value = 100000
$ws.Range("I117").Value = 100000
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 100000
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = -95411
# Row 118
$ws.Range("H118").Value = 69000
$ws.Range("I118").Value = 69000
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 69000
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -67343
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
# Row 120
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
# Row 123
$ws.Range("H123").Value = 93999
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 93999
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 93999
$ws.Range("N123").Value = -103799
# Row 124
$ws.Range("H124").Value = 67519.664
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 67519.664
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 67519.664
$ws.Range("N124").Value = -77339.664
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
# Row 129
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 0
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
# Row 131
$ws.Range("H131").Value = 92999
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 92999
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 92999
$ws.Range("N131").Value = -103079
# Row 132
$ws.Range("H132").Value = 120209
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 120209
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 120209
$ws.Range("N132").Value = -130329
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
# Row 135
$ws.Range("H135").Value = 59450
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 59450
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 59450
$ws.Range("N135").Value = -69590
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
# Row 138
$ws.Range("H138").Value = 49999
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 49999
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 49999
$ws.Range("N138").Value = -60279
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
# Row 140
$ws.Range("H140").Value = 95000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 95000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 95000
$ws.Range("N140").Value = -105360
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3000
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 2000
$ws.Range("M31").Value = -1705
# Row 34
$ws.Range("H34").Value = 3000
$ws.Range("I34").Value = 2000
$ws.Range("K34").Value = 2000
$ws.Range("M34").Value = -1798

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 132.66667
$ws.Range("I5").Value = 99
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 297
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -185
$ws.Range("N5").Value = -824
# Row 14
$ws.Range("H14").Value = 270.85715
$ws.Range("I14").Value = 270.85715
$ws.Range("K14").Value = 812.5714499999999
$ws.Range("M14").Value = -639.5714499999999
# Row 97
$ws.Range("H97").Value = 7500
$ws.Range("I97").Value = 10000
$ws.Range("J97").Value = 5000
$ws.Range("K97").Value = 30000
$ws.Range("L97").Value = 15000
$ws.Range("M97").Value = -29504
$ws.Range("N97").Value = -15992
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
# Row 113
$ws.Range("H113").Value = 413
$ws.Range("I113").Value = 449.2857
$ws.Range("J113").Value = 286
$ws.Range("K113").Value = 1347.8571
$ws.Range("L113").Value = 858
$ws.Range("M113").Value = 822.1428999999998
$ws.Range("N113").Value = -5198
# Row 135
$ws.Range("H135").Value = 132.66667
$ws.Range("I135").Value = 99
$ws.Range("J135").Value = 200
$ws.Range("K135").Value = 891
$ws.Range("L135").Value = 1800
$ws.Range("M135").Value = 1644
$ws.Range("N135").Value = -6870
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2828.5715
$ws.Range("I7").Value = 2800
$ws.Range("K7").Value = 2800
$ws.Range("M7").Value = -2688
# Row 16
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1830
$ws.Range("N16").ClearContents()
# Row 20
$ws.Range("H20").Value = 10000000
$ws.Range("J20").Value = 10000000
$ws.Range("L20").Value = 10000000
$ws.Range("N20").Value = -10000452
# Row 41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
# Row 43
$ws.Range("H43").Value = 29550
$ws.Range("I43").Value = 29100
$ws.Range("K43").Value = 29100
$ws.Range("M43").Value = -28907
# Row 126
$ws.Range("H126").Value = 2828.5715
$ws.Range("I126").Value = 2800
$ws.Range("K126").Value = 8400
$ws.Range("M126").Value = -5930
